$d = $word.ActiveDocument

# Locate the paragraph that ends the "Giving Roles to Product API" bullet list
# (the one that ends with "...give the details of products").
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Then Try to get a token from Postman*") {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not find anchor paragraph to insert after."
}

# Split a brand new (empty) paragraph in right after the anchor paragraph.
$target.Range.InsertParagraphAfter() | Out-Null
$newPara = $target.Next()

# Populate that new paragraph with the exact OOXML for the new list item
# (ListParagraph style, numId=1 bullet list, Cascadia Mono run formatting).
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="480" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Cascadia Mono" w:hAnsi="Cascadia Mono" w:cs="Cascadia Mono"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Cascadia Mono" w:hAnsi="Cascadia Mono" w:cs="Cascadia Mono"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>Simmillarly create Scope and Roles for Orders Microservice , And give the Orders.Contributor Role to the same user .Then if you will call for an Token via postman , you should get both the roles for the user(Orders.Contributor, Products.ReadOnly)</w:t></w:r></w:p>'

$newPara.Range.InsertXML($xml) | Out-Null
